# fix(publipostage): Refactor synthetic array /3
#
# The "statut" column (A) holds colored-square emoji used as a status
# marker, and "statut_label" (B) holds the matching French color name.
# This swaps the emoji family from colored squares to colored-book icons,
# and renames the "noir" (black) status to "bleu" (blue) to match the new
# book icon (⬛ -> 📘 "bleu" instead of "noir"). rouge/orange/vert keep
# their existing labels since only their icon glyph changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1, 1).End(4).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $statutCell = $ws.Cells.Item($r, 1)
    $statut = $statutCell.Value2

    if ($statut -eq "🟥") {
        $statutCell.Value = "📕"
    }
    elseif ($statut -eq "⬛") {
        $statutCell.Value = "📘"
    }
    elseif ($statut -eq "🟧") {
        $statutCell.Value = "📙"
    }
    elseif ($statut -eq "🟩") {
        $statutCell.Value = "📗"
    }

    $labelCell = $ws.Cells.Item($r, 2)
    if ($labelCell.Value2 -eq "noir") {
        $labelCell.Value = "bleu"
    }
}
